$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 79

$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 1).Value = "2020-08-17"
$ws.Cells.Item($row, 1).Style = "Normal"
$ws.Cells.Item($row, 2).Value = 525733
$ws.Cells.Item($row, 3).Value = 577531
$ws.Cells.Item($row, 4).Value = 78431
$ws.Cells.Item($row, 5).Value = 57023
$ws.Cells.Item($row, 6).Value = 26.17
